# Update the three-digit ÷ one-digit division problems/answers to the
# newly generated set, preserving all formatting (font, size, etc.)
# since Find/Replace only touches the text of the matched range.

$d = $word.ActiveDocument

$replacements = @(
    @{ old = "328÷7=46, 6";  new = "663÷8=82, 7" },
    @{ old = "724÷5=144, 4"; new = "669÷8=83, 5" },
    @{ old = "677÷3=225, 2"; new = "409÷8=51, 1" },
    @{ old = "309÷8=38, 5";  new = "559÷2=279, 1" },
    @{ old = "418÷2=209, 0"; new = "347÷4=86, 3" },
    @{ old = "660÷2=330, 0"; new = "148÷2=74, 0" },
    @{ old = "300÷4=75, 0";  new = "550÷9=61, 1" },
    @{ old = "694÷5=138, 4"; new = "710÷6=118, 2" },
    @{ old = "978÷8=122, 2"; new = "362÷4=90, 2" },
    @{ old = "744÷9=82, 6";  new = "885÷3=295, 0" },
    @{ old = "562÷7=80, 2";  new = "914÷4=228, 2" },
    @{ old = "571÷7=81, 4";  new = "470÷6=78, 2" },
    @{ old = "613÷5=122, 3"; new = "356÷3=118, 2" },
    @{ old = "391÷8=48, 7";  new = "939÷3=313, 0" },
    @{ old = "797÷2=398, 1"; new = "830÷2=415, 0" },
    @{ old = "351÷9=39, 0";  new = "303÷5=60, 3" },
    @{ old = "829÷6=138, 1"; new = "482÷3=160, 2" },
    @{ old = "661÷4=165, 1"; new = "324÷4=81, 0" },
    @{ old = "888÷4=222, 0"; new = "574÷8=71, 6" },
    @{ old = "373÷8=46, 5";  new = "316÷6=52, 4" },
    @{ old = "471÷9=52, 3";  new = "812÷8=101, 4" },
    @{ old = "414÷5=82, 4";  new = "405÷3=135, 0" },
    @{ old = "795÷7=113, 4"; new = "102÷9=11, 3" },
    @{ old = "939÷2=469, 1"; new = "687÷5=137, 2" },
    @{ old = "487÷8=60, 7";  new = "879÷3=293, 0" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
